$wb = $excel.ActiveWorkbook

# Remove the unused "Sheet2" worksheet (meta-data cleanup).
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Delete() | Out-Null

$ws = $wb.Worksheets.Item("sample")

# Record the primary-key information in the field-definition row (row 3):
# "no:I" -> "no:I:PK" and "name:T" -> "name:T:PK".
$ws.Range("A3").Value = "no:I:PK"
$ws.Range("B3").Value = "name:T:PK"

# The longer field names no longer fit the auto-fit widths, so widen the
# first two columns (and drop their "best fit" auto-sizing in the process).
$ws.Columns.Item(1).ColumnWidth = 6.6
$ws.Columns.Item(2).ColumnWidth = 12.25

# Leave the selection where the edit was made.
$ws.Range("B10").Select() | Out-Null
